# Generate Report for Handback
#
# Updates the handback-status workbook so that the "31ec876b..." row
# reflects a fresh handback/Xliff-generation pass:
#   - Overview sheet: "Latest HO Xliff Generate Date" for the
#     31ec876b-...md file is refreshed.
#   - zh-cn / de-de sheets: "Correspond Handoff Datetime" and
#     "Correspond Handback DateTime" for the 31ec876b-...md file are
#     refreshed to new timestamps.
#
# The b6188037-...md rows are not touched (their handback already
# happened earlier and is left as-is).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-07 11:07:34"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-07 11:07:29"
$wsZhCn.Range("K2").Value = "2016-09-07 11:08:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-07 11:07:34"
$wsDeDe.Range("K2").Value = "2016-09-07 11:08:21"
